$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date value for rows 2-11 from 46066 to 46070
$ws.Range("C2:C11").Value = 46070
